$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A57").Value = "purpose"
$ws.Range("B57").Value = "目的|もくてき"
$ws.Range("A58").Value = "eye"
$ws.Range("B58").Value = "目|め"
$ws.Range("A59").Value = "eye drops"
$ws.Range("B59").Value = "目薬|めぐすり"
$ws.Range("A60").Value = "the second"
$ws.Range("B60").Value = "二番目|にばんめ"
$ws.Range("A61").Value = "one's superiors"
$ws.Range("B61").Value = "目上の人|めうえのひと"
$ws.Range("A62").Value = "modern"
$ws.Range("B62").Value = "現代的|げんだいてき"
$ws.Range("A63").Value = "social"
$ws.Range("B63").Value = "社会的|しゃかいてき"
$ws.Range("A64").Value = "target"
$ws.Range("B64").Value = "的|まと"
$ws.Range("A65").Value = "Western clothes"
$ws.Range("B65").Value = "洋服|ようふく"
$ws.Range("A66").Value = "the East"
$ws.Range("B66").Value = "東洋|とうよう"
$ws.Range("A67").Value = "Western food"
$ws.Range("B67").Value = "洋食|ようしょく"
$ws.Range("A68").Value = "the Atlantic"
$ws.Range("B68").Value = "大西洋|たいせいよう"
$ws.Range("A69").Value = "clothes"
$ws.Range("B69").Value = "服|ふく"
$ws.Range("A70").Value = "Western clothes"
$ws.Range("B70").Value = "洋服|ようふく"
$ws.Range("A71").Value = "uniform"
$ws.Range("B71").Value = "制服|せいふく"
$ws.Range("A72").Value = "Japanese clothes"
$ws.Range("B72").Value = "和服|わふく"
$ws.Range("A73").Value = "cafeteria"
$ws.Range("B73").Value = "食堂|しょくどう"
$ws.Range("A74").Value = "public hall"
$ws.Range("B74").Value = "公会堂|こうかいどう"
$ws.Range("A75").Value = "dignified; imposing"
$ws.Range("B75").Value = "堂々とした|どうどうとした"
$ws.Range("A76").Value = "physical labor"
$ws.Range("B76").Value = "力仕事|ちからしごと"
$ws.Range("A77").Value = "cooperation"
$ws.Range("B77").Value = "協力|きょうりょく"
$ws.Range("A78").Value = "endeavor"
$ws.Range("B78").Value = "努力|どりょく"
$ws.Range("A79").Value = "sumo wrestler"
$ws.Range("B79").Value = "力士|りきし"
$ws.Range("A80").Value = "class"
$ws.Range("B80").Value = "授業|じゅぎょう"
$ws.Range("A81").Value = "professor"
$ws.Range("B81").Value = "教授|きょうじゅ"
$ws.Range("A82").Value = "to be given"
$ws.Range("B82").Value = "授かる|さずかる"
$ws.Range("A83").Value = "class"
$ws.Range("B83").Value = "授業|じゅぎょう"
$ws.Range("A84").Value = "occupation"
$ws.Range("B84").Value = "職業|しょくぎょう"
$ws.Range("A85").Value = "industry"
$ws.Range("B85").Value = "産業|さんぎょう"
$ws.Range("A86").Value = "service industry"
$ws.Range("B86").Value = "サービス業|サービスぎょう"
$ws.Range("A87").Value = "exam"
$ws.Range("B87").Value = "試験|しけん"
$ws.Range("A88").Value = "game; match"
$ws.Range("B88").Value = "試合|しあい"
$ws.Range("A89").Value = "entrance exam"
$ws.Range("B89").Value = "入試|にゅうし"
$ws.Range("A90").Value = "to try"
$ws.Range("B90").Value = "試みる|こころみる"
$ws.Range("A91").Value = "experiment"
$ws.Range("B91").Value = "実験|じっけん"
$ws.Range("A92").Value = "experience"
$ws.Range("B92").Value = "経験|けいけん"
$ws.Range("A93").Value = "taking examination"
$ws.Range("B93").Value = "受験|じゅけん"
$ws.Range("A94").Value = "to lend"
$ws.Range("B94").Value = "貸す|かす"
$ws.Range("A95").Value = "lending"
$ws.Range("B95").Value = "貸し出し|かしだし"
$ws.Range("A96").Value = "rental condo"
$ws.Range("B96").Value = "賃貸マンション|ちんたいマンション"
$ws.Range("A97").Value = "library"
$ws.Range("B97").Value = "図書館|としょかん"
$ws.Range("A98").Value = "map"
$ws.Range("B98").Value = "地図|ちず"
$ws.Range("A99").Value = "figure"
$ws.Range("B99").Value = "図|ず"
$ws.Range("A100").Value = "signal"
$ws.Range("B100").Value = "合図|あいず"
$ws.Range("A101").Value = "to attempt"
$ws.Range("B101").Value = "図る|はかる"
$ws.Range("A102").Value = "Japanese inn"
$ws.Range("B102").Value = "旅館|りょかん"
$ws.Range("A103").Value = "movie theater"
$ws.Range("B103").Value = "映画館|えいがかん"
$ws.Range("A104").Value = "embassy"
$ws.Range("B104").Value = "大使館|たいしかん"
$ws.Range("A105").Value = "to come to an end"
$ws.Range("B105").Value = "終わる|おわる"
$ws.Range("A106").Value = "end"
$ws.Range("B106").Value = "終わり|おわり"
$ws.Range("A107").Value = "last stop"
$ws.Range("B107").Value = "終点|しゅうてん"
$ws.Range("A108").Value = "the last..."
$ws.Range("B108").Value = "最終～|さいしゅう～"
$ws.Range("A109").Value = "homework"
$ws.Range("B109").Value = "宿題|しゅくだい"
$ws.Range("A110").Value = "boarding house"
$ws.Range("B110").Value = "下宿|げしゅく"
$ws.Range("A111").Value = "lodging"
$ws.Range("B111").Value = "宿泊|しゅくはく"
$ws.Range("A112").Value = "inn"
$ws.Range("B112").Value = "宿|やど"
$ws.Range("A113").Value = "problem; question"
$ws.Range("B113").Value = "問題|もんだい"
$ws.Range("A114").Value = "topic of conversation"
$ws.Range("B114").Value = "話題|わだい"
$ws.Range("A115").Value = "title"
$ws.Range("B115").Value = "題|だい"

Write-Output "Added rows 57-115"
